# Scheduled cryptos data refresh (GitHub Actions) - update Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.438.19"
$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").Value = "3.151.77"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.60"
$ws.Range("E5").Value = "  +0.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.07"
$ws.Range("E6").Value = "  -1.38%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "3.149.70"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.151"
$ws.Range("E10").Value = "  -0.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.37"
$ws.Range("E11").Value = "  -2.54%  "

$ws.Range("E12").Value = "  -0.78%  "

$ws.Range("E13").Value = "  +0.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.52"
$ws.Range("E14").Value = "  -1.49%  "

$ws.Range("D15").Value = "3.671.17"
$ws.Range("E15").Value = "  +0.03%  "

$ws.Range("E16").Value = "  +2.93%  "

$ws.Range("D17").Value = "64.417.44"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("D18").Value = "3.155.46"
$ws.Range("E18").Value = "  +0.16%  "

$ws.Range("E19").Value = "  -1.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.88"
$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.69"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.725"
$ws.Range("E22").Value = "  +1.77%  "

$ws.Range("E23").Value = "  +1.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.78"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.71"
$ws.Range("E25").Value = "  +1.95%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.57"
$ws.Range("E28").Value = "  +1.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.39"
$ws.Range("E29").Value = "  +8.18%  "

$ws.Range("E30").Value = "  +2.50%  "

$ws.Range("E31").Value = "  -4.25%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.49"
$ws.Range("E33").Value = "  +1.11%  "

$ws.Range("E34").Value = "  -3.72%  "

$ws.Range("E35").Value = "  +1.36%  "

$ws.Range("E36").Value = "  -0.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.79"
$ws.Range("E37").Value = "  -2.47%  "

$ws.Range("D38").Value = "0.0₃0752"
$ws.Range("E38").Value = "  +3.88%  "

$ws.Range("E39").Value = "  +5.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "453.30"
$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0395"
$ws.Range("E41").Value = "  -0.30%  "

$ws.Range("E42").Value = "  -1.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.33"
$ws.Range("E43").Value = "  -1.52%  "

$ws.Range("D44").Value = "2.844.13"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.268"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("E46").Value = "  +0.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.44"
$ws.Range("E47").Value = "  +6.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.43"
$ws.Range("E48").Value = "  +0.22%  "

$ws.Range("E49").Value = "  +0.03%  "

$ws.Range("E50").Value = "  -0.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.86"
$ws.Range("E51").Value = "  +1.97%  "
